$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.193.93"
$ws.Range("E2").Value = "'  +1.01%  "
$ws.Range("D3").Value = "'1.869.13"
$ws.Range("E3").Value = "'  +3.18%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "'  +0.33%  "
$ws.Range("D5").Value = "'311.86"
$ws.Range("E6").Value = "'  +0.22%  "
$ws.Range("D7").Value = "'0.5037"
$ws.Range("E7").Value = "'  +0.82%  "
$ws.Range("D8").Value = "'0.3903"
$ws.Range("E8").Value = "'  +0.18%  "
$ws.Range("D9").Value = "'0.09495"
$ws.Range("E9").Value = "'  -0.41%  "
$ws.Range("D10").Value = "'1.139"
$ws.Range("E10").Value = "'  +3.66%  "
$ws.Range("D11").Value = "'40.81"
$ws.Range("E11").Value = "'  +1.12%  "
$ws.Range("E12").Value = "'  +0.56%  "
$ws.Range("D13").Value = "'20.93"
$ws.Range("E13").Value = "'  +1.95%  "
$ws.Range("D14").Value = "'1.873.17"
$ws.Range("E14").Value = "'  +3.21%  "
$ws.Range("D15").Value = "'1.004"
$ws.Range("E15").Value = "'  +0.37%  "
$ws.Range("D16").Value = "'7.381"
$ws.Range("E16").Value = "'  +1.43%  "
$ws.Range("D17").Value = "'0.00001125"
$ws.Range("E17").Value = "'  -0.16%  "
$ws.Range("D18").Value = "'92.61"
$ws.Range("E18").Value = "'  -0.88%  "
$ws.Range("D19").Value = "'0.06601"
$ws.Range("E19").Value = "'  +0.10%  "
$ws.Range("D20").Value = "'17.67"
$ws.Range("E20").Value = "'  +3.00%  "
$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "'  +0.23%  "
$ws.Range("D22").Value = "'6.177"
$ws.Range("E22").Value = "'  +3.88%  "
$ws.Range("D23").Value = "'28.263.72"
$ws.Range("E23").Value = "'  +1.04%  "
$ws.Range("D24").Value = "'11.26"
$ws.Range("E24").Value = "'  +0.36%  "
$ws.Range("D25").Value = "'2.290"
$ws.Range("E25").Value = "'  +1.93%  "
$ws.Range("B26").Value = "'LidoDAOToken"
$ws.Range("C26").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.569"
$ws.Range("E26").Value = "'  +7.04%  "
$ws.Range("B27").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "'2.087.43"
$ws.Range("E27").Value = "'  +3.20%  "
$ws.Range("B28").Value = "'EthereumClassic"
$ws.Range("C28").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'21.14"
$ws.Range("E28").Value = "'  +1.78%  "
$ws.Range("B29").Value = "'Monero"
$ws.Range("C29").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'159.00"
$ws.Range("E29").Value = "'  +0.78%  "
$ws.Range("B30").Value = "'BitcoinCash"
$ws.Range("C30").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'127.26"
$ws.Range("E30").Value = "'  -0.94%  "
$ws.Range("B31").Value = "'Stellar"
$ws.Range("C31").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.1062"
$ws.Range("E31").Value = "'  -1.04%  "
$ws.Range("B32").Value = "'ImmutableX"
$ws.Range("C32").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.064"
$ws.Range("E32").Value = "'  +0.86%  "
$ws.Range("B33").Value = "'Filecoin"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.618"
$ws.Range("E33").Value = "'  +0.23%  "
$ws.Range("B34").Value = "'HuobiToken"
$ws.Range("C34").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.623"
$ws.Range("E34").Value = "'  +0.02%  "
$ws.Range("B35").Value = "'Hedera"
$ws.Range("C35").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.06740"
$ws.Range("E35").Value = "'  -1.23%  "
$ws.Range("B36").Value = "'FraxShare"
$ws.Range("C36").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'9.477"
$ws.Range("E36").Value = "'  +5.23%  "
$ws.Range("B37").Value = "'VeChain"
$ws.Range("C37").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02406"
$ws.Range("E37").Value = "'  +3.74%  "
$ws.Range("B38").Value = "'Algorand"
$ws.Range("C38").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2187"
$ws.Range("E38").Value = "'  +1.74%  "
$ws.Range("B39").Value = "'Aptos"
$ws.Range("C39").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'11.50"
$ws.Range("E39").Value = "'  +1.08%  "
$ws.Range("B40").Value = "'TheSandbox"
$ws.Range("C40").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6348"
$ws.Range("E40").Value = "'  +1.67%  "
$ws.Range("D41").Value = "'4.989"
$ws.Range("B42").Value = "'TrustWalletToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.184"
$ws.Range("E42").Value = "'  +3.35%  "
$ws.Range("B43").Value = "'Frax"
$ws.Range("C43").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "'1.003"
$ws.Range("E43").Value = "'  +0.31%  "
$ws.Range("B44").Value = "'EnergySwap"
$ws.Range("C44").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'13.44"
$ws.Range("E44").Value = "'  +2.66%  "
$ws.Range("B45").Value = "'Decentraland"
$ws.Range("C45").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.5981"
$ws.Range("E45").Value = "'  +1.03%  "
$ws.Range("B46").Value = "'WEMIXTOKEN"
$ws.Range("C46").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.279"
$ws.Range("E46").Value = "'  -1.07%  "
$ws.Range("B47").Value = "'PancakeSwap"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.660"
$ws.Range("E47").Value = "'  -0.85%  "
$ws.Range("B48").Value = "'NEARProtocol"
$ws.Range("C48").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.998"
$ws.Range("E48").Value = "'  +1.82%  "
$ws.Range("B49").Value = "'Quant"
$ws.Range("C49").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'123.26"
$ws.Range("E49").Value = "'  -0.59%  "
$ws.Range("B50").Value = "'EOS"
$ws.Range("C50").Value = "'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'1.196"
$ws.Range("E50").Value = "'  +1.48%  "
$ws.Range("B51").Value = "'Cronos"
$ws.Range("C51").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.06851"
$ws.Range("E51").Value = "'  +0.93%  "
